# Update the "Förändrad" (Changed) date column (C) for rows 2-15
# from 45174 (2023-09-05) to 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
